# Weekly refresh of "Fruta / hortaliza" data.
# The values in columns D, J, K, L, M, P (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) get reshuffled
# across the data rows (2-16) while the rest of each row (Mercado,
# Region, Categoria, Variedad, Calidad, Unidad, Origen, Kg o Unidades,
# Clasificacion) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a block, keyed by row.
$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the original values before any writes, so the shuffle can be
# applied without clobbering source data that hasn't been read yet.
$original = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Destination row -> source row (which original row's D/J/K/L/M/P values
# should now occupy this row).
$map = @{
    2  = 12
    3  = 13
    4  = 9
    5  = 11
    6  = 4
    7  = 10
    8  = 14
    9  = 2
    10 = 15
    11 = 16
    12 = 8
    13 = 7
    14 = 6
    15 = 3
    16 = 5
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $original[$srcRow][$c]
    }
}
